# add validation check before vlan creation - cisco & junos
#
# Updates the sample row in the Vlan_params sheet:
#   vlan_name : auto_test -> ashertest
#   vlan_type : users     -> Servers
#   is_dhcp   : True      -> False

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vlan_params")

$ws.Range("B2").Value = "ashertest"   # vlan_name
$ws.Range("B5").Value = "Servers"     # vlan_type

# B10 (is_dhcp) needs the literal text "False", not the Boolean value FALSE.
# A direct $ws.Range("B10").Value = "False" assignment gets auto-coerced by
# Excel into a Logical/Boolean cell (like typing FALSE straight into a
# cell would). Building the text via a formula and copying it across as a
# value keeps it as plain text, matching how the sheet was authored.
$helper = $ws.Range("Z1")
$helper.Formula = '="Fal" & "se"'
$helper.Copy()
$ws.Range("B10").PasteSpecial(-4163)  # xlPasteValues
$helper.ClearContents()

$ws.Range("B12").Select()
